$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A so the longer accessory names are readable
$ws.Columns.Item(1).ColumnWidth = 19.65

# New row 91 - Ishtar Maid Ver accessory
$ws.Cells.Item(91, 1).Value = "Ishtar Maid Ver"
$ws.Cells.Item(91, 2).Value = 5
$ws.Cells.Item(91, 3).Value = "M"
$ws.Cells.Item(91, 4).Value = 90
$ws.Cells.Item(91, 5).Value = 70
$ws.Cells.Item(91, 6).Value = 80000

# New row 92 - Kotori Police accessory
$ws.Cells.Item(92, 1).Value = "Kotori Police"
$ws.Cells.Item(92, 2).Value = 2
$ws.Cells.Item(92, 3).Value = "M"
$ws.Cells.Item(92, 4).Value = 90
$ws.Cells.Item(92, 5).Value = 70
$ws.Cells.Item(92, 6).Value = 70000

# Scroll the view down near the new rows and select B78, mirroring the
# author's final on-screen position when the rows were added
$ws.Range("B78").Select()
$excel.ActiveWindow.ScrollRow = 62
